# Weekly data refresh: a new week's record is inserted at the top of the
# data block (row 9), pushing the previously-existing rows 9-18 down to
# rows 10-19 (all of their original values are preserved verbatim).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 9 - this shifts rows 9:18 down to 10:19,
# carrying all of their existing values/formatting with them untouched.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with this week's record.
$ws.Range("A9").Value = 11
$ws.Range("B9").Value = "Vega Monumental Concepción"
$ws.Range("C9").Value = "Bíobío"
$ws.Range("D9").Value = 45233
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100107
$ws.Range("H9").Value = "Otros"
$ws.Range("I9").Value = 100107011
$ws.Range("J9").Value = "Tuna"
$ws.Range("K9").Value = "Sin especificar"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 26000
$ws.Range("O9").Value = 26000
$ws.Range("P9").Value = 26000
$ws.Range("Q9").Value = "$/caja 18 kilos"
$ws.Range("R9").Value = "Provincia de Melipilla"
$ws.Range("S9").Value = 1444
$ws.Range("T9").Value = 18

# Match the date cell styling used by the rest of the "Fecha" column (D).
$ws.Range("D9").NumberFormat = $ws.Range("D10").NumberFormat
